$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the styling of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I and J for rows 2-13
$values = @(
    @(1, 1),
    @(1, 5),
    @(1, 4),
    @(1, 6),
    @(1, 7),
    @(1, 4),
    @(1, 7),
    @(9, 10),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(4, 5)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
